# Append the latest EUR->ARS quote as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 99

# Column A holds a date-looking string ("2025-10-25"). The sheet stores every
# value as plain text (see existing rows), so prefix with a leading
# apostrophe to force text entry and avoid Excel auto-converting it to a
# date serial number.
$ws.Range("A" + $newRow).Value = "'2025-10-25"
$ws.Range("B" + $newRow).Value = "15:22:06"
$ws.Range("C" + $newRow).Value = "1.00 EUR = 1,797.6754"
